$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 110.82353
$ws.Range("I55").Value = 99.666664
$ws.Range("J55").Value = 116.90909
$ws.Range("K55").Value = 99.666664
$ws.Range("L55").Value = 116.90909
$ws.Range("M55").Value = 114.333336
$ws.Range("N55").Value = -544.90909
$ws.Range("H74").Value = 3938.3333
$ws.Range("I74").Value = 3915.9092
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3915.9092
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2979.9092
$ws.Range("N74").Value = -5872
$ws.Range("H76").Value = 3353.8462
$ws.Range("I76").Value = 3475
$ws.Range("K76").Value = 3475
$ws.Range("M76").Value = -3160
$ws.Range("H77").Value = 3938.3333
$ws.Range("I77").Value = 3915.9092
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 19579.546
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -14899.546
$ws.Range("N77").Value = -29360
$ws.Range("H79").Value = 3353.8462
$ws.Range("I79").Value = 3475
$ws.Range("K79").Value = 3475
$ws.Range("M79").Value = -2383
$ws.Range("H80").Value = 556.0345
$ws.Range("I80").Value = 403.5
$ws.Range("J80").Value = 805.63635
$ws.Range("K80").Value = 1210.5
$ws.Range("L80").Value = 2416.90905
$ws.Range("M80").Value = -212.5
$ws.Range("N80").Value = -4412.90905
$ws.Range("H83").Value = 556.0345
$ws.Range("I83").Value = 403.5
$ws.Range("J83").Value = 805.63635
$ws.Range("K83").Value = 3631.5
$ws.Range("L83").Value = 7250.72715
$ws.Range("M83").Value = 1360.5
$ws.Range("N83").Value = -17234.72715
$ws.Range("H129").Value = 994.65625
$ws.Range("J129").Value = 1019.13336
$ws.Range("L129").Value = 3057.40008
$ws.Range("N129").Value = -13057.40008
$ws.Range("H140").Value = 64990.477
$ws.Range("J140").Value = 64990.477
$ws.Range("L140").Value = 64990.477
$ws.Range("N140").Value = -75350.477

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2869.92
$ws.Range("I32").Value = 2694.8132
$ws.Range("J32").Value = 3395.24
$ws.Range("K32").Value = 2694.8132
$ws.Range("L32").Value = 3395.24
$ws.Range("M32").Value = -2407.8132
$ws.Range("N32").Value = -3969.24
$ws.Range("H61").Value = 41750884
$ws.Range("I61").Value = 45500932
$ws.Range("J61").Value = 500350
$ws.Range("K61").Value = 45500932
$ws.Range("L61").Value = 500350
$ws.Range("M61").Value = -45500720
$ws.Range("N61").Value = -500774
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -24992
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H111").Value = 45000
$ws.Range("J111").Value = 45000
$ws.Range("L111").Value = 45000
$ws.Range("N111").Value = -53180
$ws.Range("H112").Value = 26397
$ws.Range("J112").Value = 26397
$ws.Range("L112").Value = 26397
$ws.Range("N112").Value = -29351
$ws.Range("H121").Value = 64000
$ws.Range("J121").Value = 64000
$ws.Range("L121").Value = 64000
$ws.Range("N121").Value = -67494
$ws.Range("H124").Value = 31214.5
$ws.Range("J124").Value = 31214.5
$ws.Range("L124").Value = 31214.5
$ws.Range("N124").Value = -41034.5
$ws.Range("H125").Value = 55058.824
$ws.Range("J125").Value = 55058.824
$ws.Range("L125").Value = 55058.824
$ws.Range("N125").Value = -64898.824
$ws.Range("H132").Value = 66649.44
$ws.Range("I132").Value = 42473.75
$ws.Range("J132").Value = 139176.5
$ws.Range("K132").Value = 127421.25
$ws.Range("L132").Value = 417529.5
$ws.Range("M132").Value = -124891.25
$ws.Range("N132").Value = -422589.5
$ws.Range("H135").Value = 47274.08
$ws.Range("J135").Value = 47274.08
$ws.Range("L135").Value = 47274.08
$ws.Range("N135").Value = -57414.08
$ws.Range("H136").Value = 41750884
$ws.Range("I136").Value = 45500932
$ws.Range("J136").Value = 500350
$ws.Range("K136").Value = 136502796
$ws.Range("L136").Value = 1501050
$ws.Range("M136").Value = -136500246
$ws.Range("N136").Value = -1506150
$ws.Range("H139").Value = 40602.5
$ws.Range("J139").Value = 40602.5
$ws.Range("L139").Value = 40602.5
$ws.Range("N139").Value = -50882.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1465.8889
$ws.Range("I20").Value = 1281
$ws.Range("J20").Value = 1613.8
$ws.Range("K20").Value = 1281
$ws.Range("L20").Value = 1613.8
$ws.Range("M20").Value = -1034
$ws.Range("N20").Value = -2107.8
$ws.Range("H86").Value = 16700
$ws.Range("I86").Value = 20933.334
$ws.Range("J86").Value = 4000
$ws.Range("K86").Value = 20933.334
$ws.Range("L86").Value = 4000
$ws.Range("M86").Value = -19810.334
$ws.Range("N86").Value = -6246
$ws.Range("H89").Value = 16700
$ws.Range("I89").Value = 20933.334
$ws.Range("J89").Value = 4000
$ws.Range("K89").Value = 104666.67
$ws.Range("L89").Value = 20000
$ws.Range("M89").Value = -99050.67
$ws.Range("N89").Value = -31232
$ws.Range("H94").Value = 959.75
$ws.Range("I94").Value = 954
$ws.Range("K94").Value = 954
$ws.Range("M94").Value = -503

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 23500000
$ws.Range("I6").Value = 30333334
$ws.Range("J6").Value = 3000000
$ws.Range("K6").Value = 30333334
$ws.Range("L6").Value = 3000000
$ws.Range("M6").Value = -30333221
$ws.Range("N6").Value = -3000226
$ws.Range("H31").Value = 2472.5576
$ws.Range("I31").Value = 1680.2069
$ws.Range("J31").Value = 3471.6086
$ws.Range("K31").Value = 1680.2069
$ws.Range("L31").Value = 3471.6086
$ws.Range("M31").Value = -1385.2069
$ws.Range("N31").Value = -4061.6086
$ws.Range("H34").Value = 2472.5576
$ws.Range("I34").Value = 1680.2069
$ws.Range("J34").Value = 3471.6086
$ws.Range("K34").Value = 1680.2069
$ws.Range("L34").Value = 3471.6086
$ws.Range("M34").Value = -1478.2069
$ws.Range("N34").Value = -3875.6086
$ws.Range("H58").Value = 47621210
$ws.Range("I58").Value = 71430280
$ws.Range("J58").Value = 3057.2856
$ws.Range("K58").Value = 71430280
$ws.Range("L58").Value = 3057.2856
$ws.Range("M58").Value = -71430077
$ws.Range("N58").Value = -3463.2856
$ws.Range("H62").Value = 3666.6667
$ws.Range("H65").Value = 3666.6667
$ws.Range("H134").Value = 31818.236
$ws.Range("I134").Value = 2339.3
$ws.Range("J134").Value = 142364.25
$ws.Range("K134").Value = 7017.900000000001
$ws.Range("L134").Value = 427092.75
$ws.Range("M134").Value = -4482.900000000001
$ws.Range("N134").Value = -432162.75
$ws.Range("H136").Value = 47621210
$ws.Range("I136").Value = 71430280
$ws.Range("J136").Value = 3057.2856
$ws.Range("K136").Value = 214290840
$ws.Range("L136").Value = 9171.856800000001
$ws.Range("M136").Value = -214288290
$ws.Range("N136").Value = -14271.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4764436
$ws.Range("J4").Value = 4764436
$ws.Range("L4").Value = 14293308
$ws.Range("N4").Value = -14293532
$ws.Range("H6").Value = 499.5
$ws.Range("I6").Value = 49.25
$ws.Range("K6").Value = 147.75
$ws.Range("M6").Value = -34.75
$ws.Range("H12").Value = 23255898
$ws.Range("I12").Value = 52631644
$ws.Range("J12").Value = 98.208336
$ws.Range("K12").Value = 157894932
$ws.Range("L12").Value = 294.625008
$ws.Range("M12").Value = -157894759
$ws.Range("N12").Value = -640.625008
$ws.Range("H131").Value = 1372.683
$ws.Range("I131").Value = 800
$ws.Range("K131").Value = 2400
$ws.Range("M131").Value = 2640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 3582.75
$ws.Range("I31").Value = 3582.75
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 3582.75
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -3290.75
$ws.Range("N31").ClearContents()
$ws.Range("H37").Value = 3582.75
$ws.Range("I37").Value = 3582.75
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 3582.75
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -3305.75
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 807.5
$ws.Range("I46").Value = 761.6667
$ws.Range("J46").Value = 945
$ws.Range("K46").Value = 761.6667
$ws.Range("L46").Value = 945
$ws.Range("M46").Value = -573.6667
$ws.Range("N46").Value = -1321
$ws.Range("H93").Value = 2083
$ws.Range("I93").Value = 1999.5
$ws.Range("K93").Value = 1999.5
$ws.Range("M93").Value = -751.5
$ws.Range("H94").Value = 33000
$ws.Range("J94").Value = 33000
$ws.Range("L94").Value = 33000
$ws.Range("N94").Value = -34352
$ws.Range("H133").Value = 38907.41
$ws.Range("J133").Value = 39569.668
$ws.Range("L133").Value = 39569.668
$ws.Range("N133").Value = -44629.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1436.091
$ws.Range("I126").Value = 1379.4
$ws.Range("J126").Value = 1483.3334
$ws.Range("K126").Value = 4138.200000000001
$ws.Range("L126").Value = 4450.0002
$ws.Range("M126").Value = -1668.200000000001
$ws.Range("N126").Value = -9390.0002
